$d = $word.ActiveDocument

# --- Paragraph 1: the title "C. W. Editors Arrested In Air Raid Drill" ---
# The pandoc title-block export uses the "Title" style (instead of
# "Heading1") and tokenizes the heading text into one run per word /
# punctuation mark. We rebuild the paragraph text token-by-token so each
# piece is inserted as its own run.
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Style = "Title"
# Exclude the trailing paragraph mark from the working range so that
# Collapse/InsertAfter never spills over into the following paragraph.
[void]$titleRange.MoveEnd(1, -1)

$titleTokens = @("C", ".", " ", "W", ".", " ", "Editors", " ", "Arrested", " ", "In", " ", "Air", " ", "Raid", " ", "Drill")

$titleRange.Text = $titleTokens[0]
for ($i = 1; $i -lt $titleTokens.Count; $i++) {
    $titleRange.Collapse(0)
    $titleRange.InsertAfter($titleTokens[$i])
}

# --- Paragraph 2: the byline "By Dorothy Day" ---
# Pandoc emits the author list under the "Authors" style, with the leading
# "By " dropped, the bold direct-formatting removed, and the remaining
# "Dorothy Day" tokenized the same way as the title. The old paragraph
# carries hard "Bold" direct formatting, which Word has no way to fully
# strip back off once applied (re-assigning Bold=False just writes an
# explicit "off" override) -- so instead we delete the old paragraph
# outright and grow a brand new, formatting-free one in its place.
$authorPara = $d.Paragraphs.Item(2)
$authorPara.Range.Delete()

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$authorPara = $d.Paragraphs.Item(2)
$authorRange = $authorPara.Range
$authorRange.Style = "Authors"
[void]$authorRange.MoveEnd(1, -1)

$authorTokens = @("Dorothy", " ", "Day")

$authorRange.Text = $authorTokens[0]
for ($i = 1; $i -lt $authorTokens.Count; $i++) {
    $authorRange.Collapse(0)
    $authorRange.InsertAfter($authorTokens[$i])
}
